# Highlight "tipe Alpha dan Beta testing" in yellow, leaving the
# preceding " dengan " un-highlighted. This mirrors the diff where the
# run " dengan tipe " was split into " dengan " (unchanged) and "tipe "
# (highlighted), and the following runs ("Alpha", " dan ", "Beta testing")
# each gained a yellow highlight.

$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("tipe Alpha dan Beta testing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # wdYellow = 7
    $rng.Font.HighlightColorIndex = 7
}
